$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Style = "Normal"
$ws.Range("D2").Value = "'27.803.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.34%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Style = "Normal"
$ws.Range("D3").Value = "'1.905.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.51%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Style = "Normal"
$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D6").Style = "Normal"
$ws.Range("D6").Value = "'0.9996"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Style = "Normal"
$ws.Range("D7").Value = "'0.5240"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("E7").Value = "'  +7.26%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Style = "Normal"
$ws.Range("D8").Value = "'0.3782"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.32%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Style = "Normal"
$ws.Range("D9").Value = "'0.07235"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.30%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Style = "Normal"
$ws.Range("D10").Value = "'21.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.65%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Style = "Normal"
$ws.Range("D11").Value = "'0.9058"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.77%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Style = "Normal"
$ws.Range("D12").Value = "'0.07646"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.06%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Style = "Normal"
$ws.Range("D13").Value = "'1.920.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.27%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Style = "Normal"
$ws.Range("D14").Value = "'5.450"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.48%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Style = "Normal"
$ws.Range("D15").Value = "'92.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.70%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Style = "Normal"
$ws.Range("D16").Value = "'0.9993"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.18%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.03%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Style = "Normal"
$ws.Range("D18").Value = "'0.9995"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Style = "Normal"
$ws.Range("D19").Value = "'27.848.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.55%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Style = "Normal"
$ws.Range("D20").Value = "'14.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.29%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Style = "Normal"
$ws.Range("D21").Value = "'5.151"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.51%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Style = "Normal"
$ws.Range("D22").Value = "'2.133.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.40%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Style = "Normal"
$ws.Range("D23").Value = "'10.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.06%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Style = "Normal"
$ws.Range("D24").Value = "'6.629"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.16%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Style = "Normal"
$ws.Range("D25").Value = "'153.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.24%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Style = "Normal"
$ws.Range("D26").Value = "'1.869"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.80%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Style = "Normal"
$ws.Range("D27").Value = "'2.168"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.03%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Style = "Normal"
$ws.Range("D28").Value = "'18.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.17%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Style = "Normal"
$ws.Range("D29").Value = "'114.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.60%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Style = "Normal"
$ws.Range("D30").Value = "'4.861"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.57%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Style = "Normal"
$ws.Range("D31").Value = "'0.09035"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.42%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Style = "Normal"
$ws.Range("D32").Value = "'3.176"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.77%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Style = "Normal"
$ws.Range("D33").Value = "'4.847"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Style = "Normal"
$ws.Range("E33").Value = "'  +4.61%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.82%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Style = "Normal"
$ws.Range("D35").Value = "'0.7813"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.73%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Style = "Normal"
$ws.Range("D36").Value = "'0.02099"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.06%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Style = "Normal"
$ws.Range("D37").Value = "'2.621"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.50%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Style = "Normal"
$ws.Range("D38").Value = "'3.073"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.00%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Style = "Normal"
$ws.Range("D39").Value = "'0.5599"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").Style = "Normal"
$ws.Range("D40").Value = "'1.091"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.42%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Style = "Normal"
$ws.Range("D41").Value = "'0.05274"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.05%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Style = "Normal"
$ws.Range("D42").Value = "'6.733"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.25%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Style = "Normal"
$ws.Range("D43").Value = "'114.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.32%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Style = "Normal"
$ws.Range("D44").Value = "'8.567"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.25%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Style = "Normal"
$ws.Range("D45").Value = "'0.1517"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.19%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Style = "Normal"
$ws.Range("D46").Value = "'0.4805"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.38%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Style = "Normal"
$ws.Range("D47").Value = "'10.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.98%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Style = "Normal"
$ws.Range("D48").Value = "'0.9994"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.10%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Style = "Normal"
$ws.Range("D49").Value = "'1.622"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.15%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Style = "Normal"
$ws.Range("D50").Value = "'66.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.85%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Style = "Normal"
$ws.Range("D51").Value = "'0.05987"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.03%  "
$ws.Range("E51").Style = "Normal"
